$wb = $excel.ActiveWorkbook

# Sheet 2 = 建物 (building): I2 holds property_category, currently "land" -> should be "building"
$wsBuilding = $wb.Worksheets.Item(2)
$wsBuilding.Range("I2").Value = "building"

# Sheet 3 = 汽車 (car): H2 holds property_category, currently "land" -> should be "car"
$wsCar = $wb.Worksheets.Item(3)
$wsCar.Range("H2").Value = "car"
